# Apply updates described by the commit "Automatic update of files."
#  1. Column C ("Förändrad") for every existing data row (2..383) changes
#     from 45192 to 45202.
#  2. Row 383 gets an explicit row height (15, customHeight) - it previously
#     had no explicit height since it was the last row in the sheet.
#  3. Three brand-new rows (384, 385, 386) are appended with fresh
#     cleaning-notice data, all dated 45196/45202.
#  4. The sheet dimension grows from A1:Y383 to A1:Y386 (this happens
#     automatically once the new cells are populated).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Bulk-update column C (rows 2 through 383) to the new "changed" date.
$ws.Range("C2:C383").Value = 45202

# 2. Row 383 is no longer the last row, so it now carries an explicit
#    (default) row height, matching the rest of the sheet.
$ws.Rows.Item(383).RowHeight = 15

# 3. Append the three new rows.
$newRows = @(
    @{ Row = 384; A = "A 47137-2023"; G = 2.1 },
    @{ Row = 385; A = "A 47136-2023"; G = 3.7 },
    @{ Row = 386; A = "A 47134-2023"; G = 7.6 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Range("A$row").Value = $r.A

    $ws.Range("B$row").NumberFormat = "YYYY-MM-DD"
    $ws.Range("B$row").Value = 45196

    $ws.Range("C$row").NumberFormat = "YYYY-MM-DD"
    $ws.Range("C$row").Value = 45202

    $ws.Range("D$row").Value = "JÖNKÖPINGS LÄN"
    $ws.Range("E$row").Value = "VAGGERYD"

    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row").Value = 0
    $ws.Range("I$row").Value = 0
    $ws.Range("J$row").Value = 0
    $ws.Range("K$row").Value = 0
    $ws.Range("L$row").Value = 0
    $ws.Range("M$row").Value = 0
    $ws.Range("N$row").Value = 0
    $ws.Range("O$row").Value = 0
    $ws.Range("P$row").Value = 0
    $ws.Range("Q$row").Value = 0

    $ws.Range("R$row").WrapText = $true
}

# Rows 384 and 385 keep an explicit row height like the rest of the sheet;
# row 386, being the new last row, stays without one (matching the pattern
# already present in the workbook, where only the final row lacks it).
$ws.Rows.Item(384).RowHeight = 15
$ws.Rows.Item(385).RowHeight = 15
